$wb = $excel.ActiveWorkbook

# --- 4_Access and Security: add Password/User header row, trim C2 value ---
# (done before the 6_Tenants edit below so new shared strings land in the
# same order the canonical workbook uses: "Password" before "QA Tenant Test")
$access = $wb.Worksheets.Item("4_Access and Security")
$access.Range("B1:C1").Font.Bold = $true
$access.Range("B1").Value = "Password"
$access.Range("C1").Value = "User"
$access.Range("C2").Value = "system@janeirodigital.com"

# --- 6_Tenants: H2 "QA tenant test" -> "QA Tenant Test" -------------------
$tenants = $wb.Worksheets.Item("6_Tenants")
$tenants.Range("H2").NumberFormat = "@"
$tenants.Range("H2").Value = "QA Tenant Test"

# --- Selections on the various tabs (order matters: last Activate wins) ---
$access.Range("C11").Select()

$generalActivities = $wb.Worksheets.Item("1_General Activities")
$generalActivities.Range("C7").Select()

$tenants.Activate()
$tenants.Range("H8").Select()
